$d = $word.ActiveDocument

$replacements = @(
    @("789÷7=", "864÷6="),
    @("909÷2=", "252÷9="),
    @("191÷9=", "826÷9="),
    @("249÷3=", "390÷5="),
    @("753÷7=", "741÷7="),
    @("215÷8=", "532÷3="),
    @("331÷6=", "803÷8="),
    @("403÷5=", "945÷7="),
    @("247÷6=", "356÷2="),
    @("677÷4=", "930÷5="),
    @("258÷8=", "127÷9="),
    @("605÷3=", "746÷7="),
    @("937÷2=", "848÷2="),
    @("885÷5=", "736÷6="),
    @("542÷3=", "712÷4="),
    @("610÷3=", "655÷6="),
    @("903÷2=", "594÷3="),
    @("355÷9=", "351÷7="),
    @("627÷2=", "519÷7="),
    @("709÷8=", "782÷3="),
    @("942÷6=", "114÷2="),
    @("489÷4=", "537÷7="),
    @("778÷9=", "819÷4="),
    @("254÷8=", "231÷3="),
    @("234÷4=", "873÷6=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
